$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet logs one row per analysis run. This edit adds a header row and
# two more logged runs, and rewrites the original run's row with refreshed
# numbers (Await-based timing instead of Wait-based) -- so the whole A1:M4
# block is written fresh: row 1 = headers, rows 2-4 = the three runs.

# Header row (row 1)
$ws.Range("A1").Value = "Date/Time"
$ws.Range("A1").NumberFormat = "m/d/yy h:mm"
$ws.Range("B1").Value = "Method"
$ws.Range("C1").Value = "elapsedMs"
$ws.Range("D1").Value = "wordCount"
$ws.Range("E1").Value = "sentenceCount"
$ws.Range("F1").Value = "posWordCount"
$ws.Range("G1").Value = "negWordCount"
$ws.Range("H1").Value = "posWordPercentage"
$ws.Range("I1").Value = "negWordPercentage"
$ws.Range("J1").Value = "posPhraseCount"
$ws.Range("K1").Value = "negativePhraseCount"
$ws.Range("L1").Value = "posWordPercentage"
$ws.Range("M1").Value = "negPhrasePercentage"

# Row 2 (new run)
$ws.Range("A2").Value = 42585.69090277778
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"
$ws.Range("B2").Value = "Noun"
$ws.Range("C2").Value = 11364
$ws.Range("D2").Value = 8815
$ws.Range("E2").Value = 1528
$ws.Range("F2").Value = 227
$ws.Range("G2").Value = 137
$ws.Range("H2").Value = 61
$ws.Range("I2").Value = 36
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 61
$ws.Range("M2").Value = 9

# Row 3 (new run)
$ws.Range("A3").Value = 42585.694745370369
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = "Noun"
$ws.Range("C3").Value = 12883
$ws.Range("D3").Value = 8800
$ws.Range("E3").Value = 1559
$ws.Range("F3").Value = 231
$ws.Range("G3").Value = 136
$ws.Range("H3").Value = 62
$ws.Range("I3").Value = 36
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 62
$ws.Range("M3").Value = 25

# Row 4 (new run)
$ws.Range("A4").Value = 42585.698310185187
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = "Noun"
$ws.Range("C4").Value = 11563
$ws.Range("D4").Value = 8802
$ws.Range("E4").Value = 1559
$ws.Range("F4").Value = 231
$ws.Range("G4").Value = 136
$ws.Range("H4").Value = 62
$ws.Range("I4").Value = 36
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 62
$ws.Range("M4").Value = 25

# Column widths (best-fit widths, matching the authored workbook)
$ws.Columns.Item(1).ColumnWidth = 13.85546875
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Columns.Item(3).ColumnWidth = 10.5703125
$ws.Columns.Item(4).ColumnWidth = 10.85546875
$ws.Columns.Item(5).ColumnWidth = 14.5703125
$ws.Columns.Item(6).ColumnWidth = 14.42578125
$ws.Columns.Item(7).ColumnWidth = 14.5703125
$ws.Columns.Item(8).ColumnWidth = 19.28515625
$ws.Columns.Item(9).ColumnWidth = 19.42578125
$ws.Columns.Item(10).ColumnWidth = 15.5703125
$ws.Columns.Item(11).ColumnWidth = 20.28515625
$ws.Columns.Item(12).ColumnWidth = 19.28515625
$ws.Columns.Item(13).ColumnWidth = 20.5703125
